$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers need to be forced to text
# (matching the original inline-string / text cell type) by temporarily
# applying a text number format, then resetting the style back to Normal
# so no stray style index is left attached to the cell.
$textCells = @("D4", "D6", "D8", "D9", "D10", "D11", "D12", "D13", "D19", "D22", "D24", "D27", "D30", "D31", "D33", "D34", "D38", "D39", "D40", "D43", "D44", "D47", "D48", "D49", "D50", "D51")
foreach ($c in $textCells) {
    $ws.Range($c).NumberFormat = "@"
}

$ws.Range('D2').Value = '44.981.73'
$ws.Range('E2').Value = '  +4.22%  '
$ws.Range('D3').Value = '2.425.80'
$ws.Range('E3').Value = '  +2.26%  '
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('E5').Value = '  +4.76%  '
$ws.Range('D6').Value = '104.14'
$ws.Range('E6').Value = '  +8.84%  '
$ws.Range('E7').Value = '  +2.56%  '
$ws.Range('D8').Value = '0.999'
$ws.Range('E8').Value = '  -0.11%  '
$ws.Range('D9').Value = '0.531'
$ws.Range('E9').Value = '  +10.24%  '
$ws.Range('D10').Value = '35.77'
$ws.Range('E10').Value = '  +3.85%  '
$ws.Range('D11').Value = '0.0804'
$ws.Range('E11').Value = '  +1.96%  '
$ws.Range('D12').Value = '0.122'
$ws.Range('E12').Value = '  -2.66%  '
$ws.Range('D13').Value = '18.49'
$ws.Range('E13').Value = '  +1.50%  '
$ws.Range('E14').Value = '  +2.58%  '
$ws.Range('D15').Value = '2.806.11'
$ws.Range('E15').Value = '  +2.46%  '
$ws.Range('D16').Value = '2.430.08'
$ws.Range('E16').Value = '  +3.90%  '
$ws.Range('E17').Value = '  +4.26%  '
$ws.Range('D18').Value = '44.872.13'
$ws.Range('E18').Value = '  +3.95%  '
$ws.Range('D19').Value = '12.35'
$ws.Range('E19').Value = '  +3.02%  '
$ws.Range('E20').Value = '  +1.33%  '
$ws.Range('D21').Value = '0.0₃0921'
$ws.Range('D22').Value = '68.91'
$ws.Range('E22').Value = '  +1.42%  '
$ws.Range('E23').Value = '  +3.52%  '
$ws.Range('D24').Value = '2.30'
$ws.Range('E24').Value = '  +4.26%  '
$ws.Range('E25').Value = '  +2.70%  '
$ws.Range('E26').Value = '  -0.01%  '
$ws.Range('D27').Value = '25.37'
$ws.Range('E27').Value = '  +3.41%  '
$ws.Range('E28').Value = '  -7.54%  '
$ws.Range('E29').Value = '  +2.41%  '
$ws.Range('D30').Value = '33.94'
$ws.Range('E30').Value = '  +5.64%  '
$ws.Range('D31').Value = '48.88'
$ws.Range('E31').Value = '  +1.82%  '
$ws.Range('E32').Value = '  +17.05%  '
$ws.Range('D33').Value = '19.83'
$ws.Range('E33').Value = '  +11.68%  '
$ws.Range('D34').Value = '5.24'
$ws.Range('E34').Value = '  +4.22%  '
$ws.Range('E35').Value = '  +0.20%  '
$ws.Range('E36').Value = '  +4.32%  '
$ws.Range('E37').Value = '  +4.11%  '
$ws.Range('D38').Value = '4.52'
$ws.Range('E38').Value = '  +4.68%  '
$ws.Range('B39').Value = 'LidoDAOToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D39').Value = '2.87'
$ws.Range('E39').Value = '  +0.88%  '
$ws.Range('B40').Value = 'Monero'
$ws.Range('C40').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D40').Value = '125.75'
$ws.Range('E40').Value = '  -2.25%  '
$ws.Range('E41').Value = '  +2.47%  '
$ws.Range('E42').Value = '  -2.92%  '
$ws.Range('D43').Value = '21.00'
$ws.Range('E43').Value = '  -0.29%  '
$ws.Range('D44').Value = '0.0290'
$ws.Range('E44').Value = '  +4.20%  '
$ws.Range('D45').Value = '1.946.01'
$ws.Range('E45').Value = '  +0.78%  '
$ws.Range('D47').Value = '2.95'
$ws.Range('E47').Value = '  +7.87%  '
$ws.Range('D48').Value = '9.23'
$ws.Range('E48').Value = '  -0.05%  '
$ws.Range('D49').Value = '1.79'
$ws.Range('E49').Value = '  +18.00%  '
$ws.Range('D50').Value = '75.79'
$ws.Range('E50').Value = '  +6.05%  '
$ws.Range('D51').Value = '53.97'
$ws.Range('E51').Value = '  +4.71%  '

foreach ($c in $textCells) {
    $ws.Range($c).Style = "Normal"
}

Write-Host "Applied cryptos list update"